$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(12, 9).Value = "sv"
$ws.Cells.Item(12, 10).Value = "Statement-opinion"
$ws.Cells.Item(23, 9).Value = "sv"
$ws.Cells.Item(23, 10).Value = "Statement-opinion"
$ws.Cells.Item(27, 9).Value = "sv"
$ws.Cells.Item(27, 10).Value = "Statement-opinion"
$ws.Cells.Item(31, 9).Value = "ba"
$ws.Cells.Item(31, 10).Value = "Appreciation"
$ws.Cells.Item(38, 9).Value = "ba"
$ws.Cells.Item(38, 10).Value = "Appreciation"
$ws.Cells.Item(44, 9).Value = "ba"
$ws.Cells.Item(44, 10).Value = "Appreciation"
$ws.Cells.Item(47, 9).Value = "aa"
$ws.Cells.Item(47, 10).Value = "Agree/Accept"
$ws.Cells.Item(51, 9).Value = "%"
$ws.Cells.Item(51, 10).Value = "Uninterpretable"
$ws.Cells.Item(54, 9).Value = "sd"
$ws.Cells.Item(54, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(57, 9).Value = "b"
$ws.Cells.Item(57, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(60, 9).Value = "sd"
$ws.Cells.Item(60, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(63, 9).Value = "aa"
$ws.Cells.Item(63, 10).Value = "Agree/Accept"
$ws.Cells.Item(66, 9).Value = "sv"
$ws.Cells.Item(66, 10).Value = "Statement-opinion"
$ws.Cells.Item(67, 9).Value = "b"
$ws.Cells.Item(67, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(74, 9).Value = "ba"
$ws.Cells.Item(74, 10).Value = "Appreciation"
$ws.Cells.Item(93, 9).Value = "aa"
$ws.Cells.Item(93, 10).Value = "Agree/Accept"
$ws.Cells.Item(95, 9).Value = "sd"
$ws.Cells.Item(95, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(119, 9).Value = "sv"
$ws.Cells.Item(119, 10).Value = "Statement-opinion"
$ws.Cells.Item(123, 9).Value = "aa"
$ws.Cells.Item(123, 10).Value = "Agree/Accept"
$ws.Cells.Item(129, 9).Value = "sv"
$ws.Cells.Item(129, 10).Value = "Statement-opinion"
$ws.Cells.Item(130, 9).Value = "b"
$ws.Cells.Item(130, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(136, 9).Value = "sd"
$ws.Cells.Item(136, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(141, 9).Value = "b"
$ws.Cells.Item(141, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(151, 9).Value = "sv"
$ws.Cells.Item(151, 10).Value = "Statement-opinion"
$ws.Cells.Item(159, 9).Value = "aa"
$ws.Cells.Item(159, 10).Value = "Agree/Accept"
$ws.Cells.Item(176, 9).Value = "b"
$ws.Cells.Item(176, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(180, 9).Value = "sd"
$ws.Cells.Item(180, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(186, 9).Value = "b"
$ws.Cells.Item(186, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(187, 9).Value = "sd"
$ws.Cells.Item(187, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(188, 9).Value = "ba"
$ws.Cells.Item(188, 10).Value = "Appreciation"
$ws.Cells.Item(194, 9).Value = "sv"
$ws.Cells.Item(194, 10).Value = "Statement-opinion"
$ws.Cells.Item(195, 9).Value = "ba"
$ws.Cells.Item(195, 10).Value = "Appreciation"
$ws.Cells.Item(217, 9).Value = "ba"
$ws.Cells.Item(217, 10).Value = "Appreciation"
$ws.Cells.Item(230, 9).Value = "b"
$ws.Cells.Item(230, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(236, 9).Value = "b"
$ws.Cells.Item(236, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(244, 9).Value = "ba"
$ws.Cells.Item(244, 10).Value = "Appreciation"
$ws.Cells.Item(246, 9).Value = "sv"
$ws.Cells.Item(246, 10).Value = "Statement-opinion"
$ws.Cells.Item(253, 9).Value = "b"
$ws.Cells.Item(253, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(260, 9).Value = "aa"
$ws.Cells.Item(260, 10).Value = "Agree/Accept"
$ws.Cells.Item(262, 9).Value = "aa"
$ws.Cells.Item(262, 10).Value = "Agree/Accept"
$ws.Cells.Item(267, 9).Value = "sv"
$ws.Cells.Item(267, 10).Value = "Statement-opinion"
$ws.Cells.Item(272, 9).Value = "sd"
$ws.Cells.Item(272, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(281, 9).Value = "b"
$ws.Cells.Item(281, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(283, 9).Value = "sv"
$ws.Cells.Item(283, 10).Value = "Statement-opinion"
$ws.Cells.Item(285, 9).Value = "aa"
$ws.Cells.Item(285, 10).Value = "Agree/Accept"
$ws.Cells.Item(287, 9).Value = "sd"
$ws.Cells.Item(287, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(289, 9).Value = "aa"
$ws.Cells.Item(289, 10).Value = "Agree/Accept"
$ws.Cells.Item(293, 9).Value = "aa"
$ws.Cells.Item(293, 10).Value = "Agree/Accept"
$ws.Cells.Item(301, 9).Value = "b"
$ws.Cells.Item(301, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(302, 9).Value = "ba"
$ws.Cells.Item(302, 10).Value = "Appreciation"
$ws.Cells.Item(314, 9).Value = "sv"
$ws.Cells.Item(314, 10).Value = "Statement-opinion"
$ws.Cells.Item(336, 9).Value = "ba"
$ws.Cells.Item(336, 10).Value = "Appreciation"
$ws.Cells.Item(342, 9).Value = "sd"
$ws.Cells.Item(342, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(350, 9).Value = "ba"
$ws.Cells.Item(350, 10).Value = "Appreciation"
$ws.Cells.Item(357, 9).Value = "aa"
$ws.Cells.Item(357, 10).Value = "Agree/Accept"
$ws.Cells.Item(359, 9).Value = "aa"
$ws.Cells.Item(359, 10).Value = "Agree/Accept"
$ws.Cells.Item(361, 9).Value = "aa"
$ws.Cells.Item(361, 10).Value = "Agree/Accept"
$ws.Cells.Item(366, 9).Value = "aa"
$ws.Cells.Item(366, 10).Value = "Agree/Accept"
$ws.Cells.Item(380, 9).Value = "aa"
$ws.Cells.Item(380, 10).Value = "Agree/Accept"
$ws.Cells.Item(382, 9).Value = "ba"
$ws.Cells.Item(382, 10).Value = "Appreciation"
$ws.Cells.Item(394, 9).Value = "ba"
$ws.Cells.Item(394, 10).Value = "Appreciation"
$ws.Cells.Item(397, 9).Value = "ba"
$ws.Cells.Item(397, 10).Value = "Appreciation"
$ws.Cells.Item(405, 9).Value = "sd"
$ws.Cells.Item(405, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(427, 9).Value = "sd"
$ws.Cells.Item(427, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(431, 9).Value = "b"
$ws.Cells.Item(431, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(433, 9).Value = "sv"
$ws.Cells.Item(433, 10).Value = "Statement-opinion"
$ws.Cells.Item(435, 9).Value = "sd"
$ws.Cells.Item(435, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(439, 9).Value = "%"
$ws.Cells.Item(439, 10).Value = "Uninterpretable"
$ws.Cells.Item(440, 9).Value = "b"
$ws.Cells.Item(440, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(445, 9).Value = "sd"
$ws.Cells.Item(445, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(466, 9).Value = "%"
$ws.Cells.Item(466, 10).Value = "Uninterpretable"
$ws.Cells.Item(473, 9).Value = "sd"
$ws.Cells.Item(473, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(480, 9).Value = "b"
$ws.Cells.Item(480, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(481, 9).Value = "sd"
$ws.Cells.Item(481, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(482, 9).Value = "sd"
$ws.Cells.Item(482, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(483, 9).Value = "b"
$ws.Cells.Item(483, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(484, 9).Value = "b"
$ws.Cells.Item(484, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(489, 9).Value = "aa"
$ws.Cells.Item(489, 10).Value = "Agree/Accept"
$ws.Cells.Item(491, 9).Value = "aa"
$ws.Cells.Item(491, 10).Value = "Agree/Accept"
$ws.Cells.Item(495, 9).Value = "aa"
$ws.Cells.Item(495, 10).Value = "Agree/Accept"
$ws.Cells.Item(516, 9).Value = "sd"
$ws.Cells.Item(516, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(523, 9).Value = "sd"
$ws.Cells.Item(523, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(565, 9).Value = "ba"
$ws.Cells.Item(565, 10).Value = "Appreciation"
$ws.Cells.Item(570, 9).Value = "sd"
$ws.Cells.Item(570, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(579, 9).Value = "aa"
$ws.Cells.Item(579, 10).Value = "Agree/Accept"
$ws.Cells.Item(583, 9).Value = "sd"
$ws.Cells.Item(583, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(590, 9).Value = "sd"
$ws.Cells.Item(590, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(591, 9).Value = "aa"
$ws.Cells.Item(591, 10).Value = "Agree/Accept"
$ws.Cells.Item(604, 9).Value = "sd"
$ws.Cells.Item(604, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(606, 9).Value = "sd"
$ws.Cells.Item(606, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(625, 9).Value = "ba"
$ws.Cells.Item(625, 10).Value = "Appreciation"
$ws.Cells.Item(628, 9).Value = "sv"
$ws.Cells.Item(628, 10).Value = "Statement-opinion"
$ws.Cells.Item(630, 9).Value = "b"
$ws.Cells.Item(630, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(633, 9).Value = "aa"
$ws.Cells.Item(633, 10).Value = "Agree/Accept"
$ws.Cells.Item(634, 9).Value = "b"
$ws.Cells.Item(634, 10).Value = "Acknowledge (Backchannel)"
